$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: renumber and change dish name ---
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Vector"
$ws.Cells.Item(2, 3).Value = "Шашлык Утка"
$ws.Cells.Item(2, 4).Value = 1900

# --- Row 3: the original order (previously row2 text), now row 3 ---
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Vector"
$ws.Cells.Item(3, 3).Value = "Шашлык Баранина"
$ws.Cells.Item(3, 4).Value = 1900

# --- Row 4: new combined order ---
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Vector"
$ws.Cells.Item(4, 3).Value = "Шашлык Утка - Шашлык Баранина - Кока Кола2л"
$ws.Cells.Item(4, 4).Value = 5600

# --- Row 5: zero-value order ---
$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = "Vector"
$ws.Cells.Item(5, 4).Value = 0

# --- Row 6: totals row, shaded with a new fill ---
$ws.Range("A6:D6").Interior.Color = 13434828
$ws.Cells.Item(6, 4).Value = 9400

# --- Widen the "Заказ" (order) column to fit the longer combined text ---
$ws.Columns.Item(3).ColumnWidth = 44.16666666666667
